$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header values in row 1 for the two new columns (P, Q),
# carrying over the same (bold/centered/bordered) formatting as O1.
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15
$ws.Range("O1").Copy()
$ws.Range("P1:Q1").PasteSpecial(-4122)

# For each data row (2-25): swap columns I<->K and M<->O, and
# populate the two new columns P and Q with value 2.
for ($r = 2; $r -le 25; $r++) {
    $iVal = $ws.Cells.Item($r, 9).Value()   # column I
    $kVal = $ws.Cells.Item($r, 11).Value()  # column K
    $mVal = $ws.Cells.Item($r, 13).Value()  # column M
    $oVal = $ws.Cells.Item($r, 15).Value()  # column O

    $ws.Cells.Item($r, 9).Value = $kVal   # I = old K
    $ws.Cells.Item($r, 11).Value = $iVal  # K = old I
    $ws.Cells.Item($r, 13).Value = $oVal  # M = old O
    $ws.Cells.Item($r, 15).Value = $mVal  # O = old M

    $ws.Cells.Item($r, 16).Value = 2      # column P
    $ws.Cells.Item($r, 17).Value = 2      # column Q
}
